# Update country data files
#
# Adds the MSME size-classification table (new rows 18-22, with a bold
# header row 18) below the existing indicator table, and moves the
# "DGEEC" source-attribution rows from rows 23-24 down to rows 29-30 so
# they sit below the new table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Remove the existing "DGEEC" source rows (they move down to 29/30) ----
$ws.Range("A23").Clear()
$ws.Range("A24").Clear()

# ---- Header row for the new table (bold, like the other "title" rows) ----
$ws.Range("B18").Value = "Number of employees"
$ws.Range("C18").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D18").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B18:D18").Font.Bold = $true

# ---- Micro row ----
$ws.Range("A19").Value = "Micro"
$ws.Range("B19").Value = "<5"
$ws.Range("C19").Value = "< G.23 Millionlon (USD4,400 approx.)"
$ws.Range("D19").Value = "< G.70 Millionlon (USD13,4450 approx.)"

# ---- Small row ----
$ws.Range("A20").Value = "Small"
$ws.Range("B20").Value = "6-20"
$ws.Range("C20").Value = "< G.92 Millionlon (USD17,700 approx.)"
$ws.Range("D20").Value = "< G.271 Millionlon (USD52,000 approx.)"

# ---- Medium row ----
$ws.Range("A21").Value = "Medium"
$ws.Range("B21").Value = "21-100"
$ws.Range("C21").Value = "< G.460 Millionlon (USD88,500 approx.)"
$ws.Range("D21").Value = "< G.1,355 Millionlon (USD 260,600approx.)"

# ---- Large row ----
$ws.Range("A22").Value = "Large"
$ws.Range("B22").Value = ">100"
$ws.Range("C22").Value = "> G.460 Millionlon (USD88,500 approx.)"
$ws.Range("D22").Value = "> G.1,355 Millionlon (USD 260,600approx.)"

# ---- Restore the source attribution further down the sheet ----
$ws.Range("A29").Value = "DGEEC"
$ws.Range("A29").Font.Bold = $true

$ws.Range("A30").Value = 'Dirección General de Estadística, Encuestas y Censos (DGEEC), "Censo Económico Nacional 2011", 2013, p.57. Available at http://www.dgeec.gov.py/Publicaciones/Biblioteca/CEN2011/resultados_finales_CEN.pdf'
$ws.Range("A30").Font.Italic = $true

# ---- Re-apply formatting to the pre-existing rows that this headless
#      engine's load/save round trip otherwise strips (bold/italic/size/
#      underline carried by the workbook's named cell styles) ----
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = 2
$ws.Range("B9:D9").Font.Bold = $true
$ws.Range("A10").Font.Bold = $true
$ws.Range("A11").Font.Bold = $true
$ws.Range("A12").Font.Bold = $true
$ws.Range("A13").Font.Bold = $true
$ws.Range("A14").Font.Bold = $true
$ws.Range("A15").Font.Italic = $true

Write-Host "Done"
